# Scheduled-runner price refresh for the Exodus_Profits workbook.
# Updates currentAveragePrice(NQ/HQ) and the dependent Leve profit
# columns on each crafting-class sheet with freshly polled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 92494.91
$ws.Range("J109").Value = 92494.91
$ws.Range("L109").Value = 92494.91
$ws.Range("N109").Value = -95268.91
$ws.Range("H110").Value = 34323.715
$ws.Range("J110").Value = 34323.715
$ws.Range("L110").Value = 34323.715
$ws.Range("N110").Value = -42503.715
$ws.Range("H117").Value = 90027.91
$ws.Range("J117").Value = 90027.91
$ws.Range("L117").Value = 90027.91
$ws.Range("N117").Value = -99205.91

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 67779.60000000001
$ws.Range("J107").Value = 67779.60000000001
$ws.Range("L107").Value = 67779.60000000001
$ws.Range("N107").Value = -75459.60000000001
$ws.Range("H108").Value = 88854.28999999999
$ws.Range("J108").Value = 88854.28999999999
$ws.Range("L108").Value = 88854.28999999999
$ws.Range("N108").Value = -96534.28999999999
$ws.Range("H132").Value = 2235.7222
$ws.Range("I132").Value = 2174.7144
$ws.Range("K132").Value = 6524.1432
$ws.Range("M132").Value = -3994.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 54699.6
$ws.Range("J50").Value = 55697.25
$ws.Range("L50").Value = 55697.25
$ws.Range("N50").Value = -56845.25
$ws.Range("H53").Value = 68992
$ws.Range("J53").Value = 68992
$ws.Range("L53").Value = 68992
$ws.Range("N53").Value = -70140
$ws.Range("H55").Value = 35097.4
$ws.Range("J55").Value = 35097.4
$ws.Range("L55").Value = 35097.4
$ws.Range("N55").Value = -35643.4
$ws.Range("H110").Value = 80922.5
$ws.Range("J110").Value = 80922.5
$ws.Range("L110").Value = 80922.5
$ws.Range("N110").Value = -89102.5
$ws.Range("H117").Value = 96241.664
$ws.Range("J117").Value = 96241.664
$ws.Range("L117").Value = 96241.664
$ws.Range("N117").Value = -105419.664
$ws.Range("H118").Value = 71579.336
$ws.Range("H122").Value = 67968
$ws.Range("J122").Value = 67968
$ws.Range("L122").Value = 67968
$ws.Range("N122").Value = -77768
$ws.Range("H127").Value = 61478.668
$ws.Range("J127").Value = 61478.668
$ws.Range("L127").Value = 61478.668
$ws.Range("N127").Value = -71398.66800000001
$ws.Range("H132").Value = 31153.54
$ws.Range("J132").Value = 31153.54
$ws.Range("L132").Value = 31153.54
$ws.Range("N132").Value = -41273.54
$ws.Range("H134").Value = 6393.3335
$ws.Range("I134").Value = 4048.3333
$ws.Range("K134").Value = 12144.9999
$ws.Range("M134").Value = -9609.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 41698.2
$ws.Range("J9").Value = 41698.2
$ws.Range("L9").Value = 41698.2
$ws.Range("N9").Value = -42034.2
$ws.Range("H18").Value = 51944.5
$ws.Range("J18").Value = 51944.5
$ws.Range("L18").Value = 51944.5
$ws.Range("N18").Value = -52404.5
$ws.Range("H31").Value = 2404.818
$ws.Range("I31").Value = 1593.0454
$ws.Range("K31").Value = 1593.0454
$ws.Range("M31").Value = -1298.0454
$ws.Range("H34").Value = 2404.818
$ws.Range("I34").Value = 1593.0454
$ws.Range("K34").Value = 1593.0454
$ws.Range("M34").Value = -1391.0454
$ws.Range("H105").Value = 3895.7144
$ws.Range("I105").Value = 3461.6667
$ws.Range("J105").Value = 6500
$ws.Range("K105").Value = 3461.6667
$ws.Range("L105").Value = 6500
$ws.Range("M105").Value = -1714.6667
$ws.Range("N105").Value = -9994
$ws.Range("H108").Value = 75258
$ws.Range("J108").Value = 75258
$ws.Range("L108").Value = 75258
$ws.Range("N108").Value = -82938
$ws.Range("H114").Value = 39984.5
$ws.Range("J114").Value = 39984.5
$ws.Range("L114").Value = 39984.5
$ws.Range("N114").Value = -48662.5
$ws.Range("H117").Value = 38248.668
$ws.Range("J117").Value = 38248.668
$ws.Range("L117").Value = 38248.668
$ws.Range("N117").Value = -47426.668
$ws.Range("H118").Value = 72775.78
$ws.Range("J118").Value = 72775.78
$ws.Range("L118").Value = 72775.78
$ws.Range("N118").Value = -76089.78
$ws.Range("H134").Value = 1768835.1
$ws.Range("I134").Value = 2234453.8
$ws.Range("J134").Value = 113302.336
$ws.Range("K134").Value = 6703361.399999999
$ws.Range("L134").Value = 339907.008
$ws.Range("M134").Value = -6700826.399999999
$ws.Range("N134").Value = -344977.008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5680.8
$ws.Range("I56").Value = 5680.8
$ws.Range("K56").Value = 5680.8
$ws.Range("M56").Value = -5150.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10171.143
$ws.Range("J52").Value = 8399.799999999999
$ws.Range("L52").Value = 8399.799999999999
$ws.Range("N52").Value = -8917.799999999999
$ws.Range("H70").Value = 10599
$ws.Range("I70").Value = 10331.667
$ws.Range("K70").Value = 10331.667
$ws.Range("M70").Value = -10061.667
$ws.Range("H73").Value = 10599
$ws.Range("I73").Value = 10331.667
$ws.Range("K73").Value = 10331.667
$ws.Range("M73").Value = -9395.666999999999
$ws.Range("H93").Value = 19363
$ws.Range("J93").Value = 19363
$ws.Range("L93").Value = 19363
$ws.Range("N93").Value = -23107
$ws.Range("H108").Value = 50995.2
$ws.Range("J108").Value = 50995.2
$ws.Range("L108").Value = 50995.2
$ws.Range("N108").Value = -58675.2
$ws.Range("H109").Value = 33447.715
$ws.Range("J109").Value = 33447.715
$ws.Range("L109").Value = 33447.715
$ws.Range("N109").Value = -35527.715
$ws.Range("H110").Value = 99999
$ws.Range("J110").Value = 99999
$ws.Range("L110").Value = 99999
$ws.Range("N110").Value = -108179
$ws.Range("H116").Value = 59996.57
$ws.Range("J116").Value = 59996.57
$ws.Range("L116").Value = 59996.57
$ws.Range("N116").Value = -69174.57000000001
$ws.Range("H119").Value = 51084.727
$ws.Range("J119").Value = 51084.727
$ws.Range("L119").Value = 51084.727
$ws.Range("N119").Value = -60760.727
$ws.Range("H122").Value = 3848.5
$ws.Range("I122").Value = 3914.1667
$ws.Range("K122").Value = 11742.5001
$ws.Range("M122").Value = -9292.500100000001
$ws.Range("H135").Value = 46071.43
$ws.Range("J135").Value = 46071.43
$ws.Range("L135").Value = 46071.43
$ws.Range("N135").Value = -56211.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3892.9285
$ws.Range("I7").Value = 2597
$ws.Range("K7").Value = 2597
$ws.Range("M7").Value = -2485
$ws.Range("H61").Value = 1193
$ws.Range("I61").Value = 1193
$ws.Range("K61").Value = 1193
$ws.Range("M61").Value = -991
$ws.Range("H113").Value = 1193
$ws.Range("I113").Value = 1193
$ws.Range("K113").Value = 1193
$ws.Range("M113").Value = 977
$ws.Range("H118").Value = 50054.285
$ws.Range("J118").Value = 50054.285
$ws.Range("L118").Value = 50054.285
$ws.Range("N118").Value = -53368.285
$ws.Range("H121").Value = 64158.727
$ws.Range("J121").Value = 64158.727
$ws.Range("L121").Value = 64158.727
$ws.Range("N121").Value = -67652.727
$ws.Range("H126").Value = 3892.9285
$ws.Range("I126").Value = 2597
$ws.Range("K126").Value = 7791
$ws.Range("M126").Value = -5321
$ws.Range("H129").Value = 72677.8
$ws.Range("J129").Value = 71749.75
$ws.Range("L129").Value = 71749.75
$ws.Range("N129").Value = -81749.75
$ws.Range("H139").Value = 49999.5
$ws.Range("I139").Value = 49999.5
$ws.Range("K139").Value = 49999.5
$ws.Range("M139").Value = -44859.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 64499
$ws.Range("I49").Value = 40000
$ws.Range("K49").Value = 40000
$ws.Range("M49").Value = -39770
$ws.Range("H81").Value = 29592.285
$ws.Range("I81").Value = 1299
$ws.Range("K81").Value = 2598
$ws.Range("M81").Value = -1537
$ws.Range("H84").Value = 29592.285
$ws.Range("I84").Value = 1299
$ws.Range("K84").Value = 12990
$ws.Range("M84").Value = -7686
$ws.Range("H98").Value = 590
$ws.Range("J98").Value = 590
$ws.Range("L98").Value = 590
$ws.Range("N98").Value = -6580
$ws.Range("H121").Value = 36435.5
$ws.Range("J121").Value = 36435.5
$ws.Range("L121").Value = 36435.5
$ws.Range("N121").Value = -39929.5
$ws.Range("H122").Value = 1842.25
$ws.Range("J122").Value = 1826.1052
$ws.Range("L122").Value = 5478.3156
$ws.Range("N122").Value = -10378.3156
$ws.Range("H132").Value = 1310.375
$ws.Range("I132").Value = 918.3333
$ws.Range("K132").Value = 2754.9999
$ws.Range("M132").Value = -224.9998999999998
